$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 2).Value = "PRIVET_FROM_PARSER!"
}
